# Add a new worksheet "ODI Bowling Extra" right after the existing
# "ODI Batting Extra" sheet, and populate it with the scraped extra
# bowling-attribute data (mirrors the layout already used by the
# "ODI Batting Extra" sheet: a MATCH_CODE key column plus a couple of
# derived per-match metrics).

$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item("ODI Batting Extra")
$newSheet = $wb.Worksheets.Add([Type]::Missing, $srcSheet)
$newSheet.Name = "ODI Bowling Extra"

# Header row (bold, thin-bordered, centered - same look as the other
# "Extra" sheet headers in this workbook).
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "MAIDEN_OVERS"
$newSheet.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

$headerRange = $newSheet.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$matchCodes = @("4287","4297","4303","4308","4314","4321","4326","4331","4336","4342","4346","4354","4355","4429","4430","4431","4698","4700","4711","4717")

$maidenOvers = @{
    "4287" = $null
    "4297" = "0"
    "4303" = $null
    "4308" = $null
    "4314" = "2"
    "4321" = "1"
    "4326" = $null
    "4331" = "2"
    "4336" = "0"
    "4342" = "0"
    "4346" = $null
    "4354" = "0"
    "4355" = $null
    "4429" = "0"
    "4430" = "2"
    "4431" = "0"
    "4698" = "0"
    "4700" = "1"
    "4711" = "0"
    "4717" = "1"
}

$percentWickets = @{
    "4287" = $null
    "4297" = "10.00%"
    "4303" = $null
    "4308" = $null
    "4314" = "30.00%"
    "4321" = "30.00%"
    "4326" = $null
    "4331" = "30.00%"
    "4336" = "10.00%"
    "4342" = $null
    "4346" = $null
    "4354" = "20.00%"
    "4355" = $null
    "4429" = "30.00%"
    "4430" = "30.00%"
    "4431" = "10.00%"
    "4698" = "10.00%"
    "4700" = "60.00%"
    "4711" = "20.00%"
    "4717" = "30.00%"
}

# All values on this sheet (including the MATCH_CODE key) are plain text,
# matching the rest of the workbook - prefix with a leading apostrophe so
# Excel stores them as text instead of auto-converting to numbers/percentages.
$rowIndex = 2
foreach ($code in $matchCodes) {
    $newSheet.Cells.Item($rowIndex, 1).Value = "'" + $code

    $mo = $maidenOvers[$code]
    if ($null -ne $mo) {
        $newSheet.Cells.Item($rowIndex, 2).Value = "'" + $mo
    }

    $pw = $percentWickets[$code]
    if ($null -ne $pw) {
        $newSheet.Cells.Item($rowIndex, 3).Value = "'" + $pw
    }

    $rowIndex++
}

$newSheet.Columns.Item(1).ColumnWidth = 10
$newSheet.Columns.Item(2).ColumnWidth = 12
$newSheet.Columns.Item(3).ColumnWidth = 20
